# Updated symbol list on Fri Dec 16 07:21:31 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Column D (Price) numeric-looking text updates ---
Set-TextValue "D2"  "263.18"
Set-TextValue "D3"  "22.87"
Set-TextValue "D4"  "6.179"
Set-TextValue "D5"  "0.06239"
Set-TextValue "D6"  "6.728"
Set-TextValue "D7"  "3.449"
Set-TextValue "D8"  "1.343"
Set-TextValue "D9"  "0.7968"
Set-TextValue "D10" "0.1589"
Set-TextValue "D11" "0.08135"
Set-TextValue "D12" "0.03427"
Set-TextValue "D13" "0.03079"
Set-TextValue "D14" "0.09328"
Set-TextValue "D15" "3.736"
Set-TextValue "D16" "0.001679"
Set-TextValue "D17" "0.04774"
Set-TextValue "D18" "0.0006131"
Set-TextValue "D19" "0.006230"
Set-TextValue "D20" "0.006188"
Set-TextValue "D21" "0.001093"
Set-TextValue "D22" "0.0001498"
Set-TextValue "D23" "3.720"
Set-TextValue "D24" "2.211"
Set-TextValue "D26" "0.1276"
Set-TextValue "D27" "0.0003197"
Set-TextValue "D40" "0.04615"
Set-TextValue "D44" "0.01011"
Set-TextValue "D46" "0.00005876"
Set-TextValue "D48" "0.6989"
Set-TextValue "D49" "0.08777"
Set-TextValue "D50" "0.00002097"

# --- Rows 41-43: coin list rotated (Kick -> BKEX -> CEJI -> Kick) ---
# Row 41 becomes BKEXToken
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D41" "0.1126"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42 becomes CEJI
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003126"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43 becomes KickToken
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D43" "0.003353"
$ws.Range("E43").Value = "42KickTokenKICK"
